# Estado de Cuenta NIT-9002476843
# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# This adds a new mora-period row (2509) below the existing 2508 row,
# updates the aggregate totals (Valor Mora / Cant. Periodos) accordingly,
# and leaves the trailing signature block shifted down by the inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the existing data row (16) so the
# signature block (previously rows 21-22) shifts down to rows 22-23.
$ws.Rows(17).Insert()

# Duplicate the look & feel (borders/fill/font/number-format) of the
# existing data row into the newly inserted row.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new data row for period 2509 (same worker/employer,
# same monthly amounts as period 2508).
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "79051517"
$ws.Range("D17").Value = "RICARDO EDULFO QUIROGA ROJAS"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 160000
$ws.Range("G17").Value = 4000000

# Update the summary figures: two overdue periods now, so the total
# overdue amount doubles.
$ws.Range("F13").Value = 2
$ws.Range("E11").Value = 320000
